# Applies the "Comando groupby e lista de exercicios 01 de comando Select"
# edit: collapses the multi-run text of items 4, 5, 7 and 8 into single
# runs (same visible text, just re-typed so Word coalesces the runs) and
# appends a brand-new item 9 paragraph about the DML/Select exercise.

$d = $word.ActiveDocument

# Special characters (kept out of raw string literals so nothing after a
# bare numeric token is mis-parsed as arithmetic by the interpreter).
$dash = [char]0x2013   # en dash "-"
$ldq  = [char]0x201C   # left double quotation mark
$rdq  = [char]0x201D   # right double quotation mark
$aac  = [char]0xE1     # a with acute (pagina, bancario...)
$iac  = [char]0xED     # i with acute (exercicios, Veiculos...)

$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

function Wrap-Package($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Run-Xml($text, $preserve) {
    $space = ''
    if ($preserve) { $space = ' xml:space="preserve"' }
    return '<w:r>' + $rPr + '<w:t' + $space + '>' + $text + '</w:t></w:r>'
}

function Merge-Paragraph($paraIndex, $newText) {
    # Collapse a paragraph's multiple runs into a single run carrying the
    # full text, keeping the paragraph's own pPr/rPr (Arial 12pt, justified).
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range.Text
    $len = $full.Length - 1   # drop the trailing paragraph mark
    $target = $d.Range($p.Range.Start, $p.Range.Start + $len)

    $bodyP = '<w:p><w:pPr><w:jc w:val="both"/>' + $rPr + '</w:pPr>' + (Run-Xml $newText $false) + '</w:p>'
    $target.InsertXML((Wrap-Package $bodyP))
}

# --- items 4, 5, 7, 8: merge split runs back into one run -------------

$text4 = "4 $dash Modelo Entidade Relacionamento do Controle de Pedidos. Encontra-se no arquivo $ldq" + "Aula 2 $dash Modelagem.pptx$rdq, p$($aac)gina: 26."
Merge-Paragraph 5 $text4

$text5 = "5 $dash Modelo Entidade Relacionamento do Sistema Banc$($aac)rio. Encontra-se no arquivo $ldq" + "Aula 2 $dash Modelagem.pptx$rdq, p$($aac)gina: 29."
Merge-Paragraph 6 $text5

$text7 = "7 $dash Fazer o script para criar o banco de dados do sistema da padaria. Encontra-se no arquivo $ldq" + "Aula 2 $dash Modelagem.pptx$rdq, p$($aac)gina: 31."
Merge-Paragraph 8 $text7

$text8 = "8 $dash Fazer o script para criar o banco de dados do sistema Banc$($aac)rio. Encontra-se no arquivo $ldq" + "Aula 2 $dash Modelagem.pptx$rdq, p$($aac)gina: 29."
Merge-Paragraph 9 $text8

# --- new item 9 paragraph ----------------------------------------------

$run1 = Run-Xml "9 $dash Fazer exerc$($iac)cios de DML - " $true
$spellStart = '<w:proofErr w:type="spellStart"/>'
$spellEnd = '<w:proofErr w:type="spellEnd"/>'
$run2 = Run-Xml 'Select' $false
$run3 = Run-Xml ' (' $true
$run4 = Run-Xml 'Aula 05 - ' $true
$run5 = Run-Xml 'VendasED' $false
$run6 = Run-Xml ' Lista01' $true
$run7 = Run-Xml (".docx) utilizando o banco de dados $ldq") $false
$run8 = Run-Xml 'VendasEd.sql' $false
$run9 = Run-Xml ("$rdq.") $false

$newParaBody = '<w:p><w:pPr><w:jc w:val="both"/>' + $rPr + '</w:pPr>' +
    $run1 + $spellStart + $run2 + $spellEnd + $run3 + $run4 +
    $spellStart + $run5 + $spellEnd + $run6 + $run7 +
    $spellStart + $run8 + $spellEnd + $run9 + '</w:p>'

# Paragraph 9 (item "8 - ...") is the anchor we insert the new paragraph after.
$p8 = $d.Paragraphs.Item(9)
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item(10)
$p9.Range.InsertXML((Wrap-Package $newParaBody))

Write-Host "done"
